$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking strings (prices) that must remain
# plain text, matching the original inline-string cells. Force text format
# before assigning, then restore the default style so no stray formatting
# is introduced.
function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "61.422.98"
$ws.Range("E2").Value = "  -0.84%  "
Set-TextValue "D3" "3.378.12"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D5" "136.74"
$ws.Range("E5").Value = "  +11.20%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue "D6" "407.28"
$ws.Range("E6").Value = "  -1.42%  "
Set-TextValue "D7" "0.593"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("E8").Value = "  +0.02%  "
Set-TextValue "D9" "0.672"
$ws.Range("E9").Value = "  +3.39%  "
$ws.Range("E10").Value = "  -4.11%  "
Set-TextValue "D11" "43.24"
$ws.Range("E11").Value = "  +4.71%  "
$ws.Range("E12").Value = "  -1.08%  "
Set-TextValue "D13" "3.904.13"
$ws.Range("E13").Value = "  -1.73%  "
Set-TextValue "D14" "8.39"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("E15").Value = "  +0.46%  "
Set-TextValue "D16" "3.375.67"
$ws.Range("E16").Value = "  -1.63%  "
Set-TextValue "D17" "61.392.26"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("E18").Value = "  -0.73%  "
Set-TextValue "D19" "11.01"
$ws.Range("E19").Value = "  +2.32%  "
Set-TextValue "D20" "0.0000127"
$ws.Range("E20").Value = "  -4.00%  "
$ws.Range("E21").Value = "  -2.60%  "
Set-TextValue "D22" "83.36"
$ws.Range("E22").Value = "  +2.56%  "
Set-TextValue "D23" "313.70"
$ws.Range("E23").Value = "  +0.43%  "
Set-TextValue "D24" "12.88"
$ws.Range("E24").Value = "  -0.51%  "
Set-TextValue "D25" "3.15"
$ws.Range("E25").Value = "  -0.38%  "
Set-TextValue "D26" "4.78"
$ws.Range("E26").Value = "  +11.57%  "
Set-TextValue "D27" "8.31"
$ws.Range("E27").Value = "  +7.93%  "
Set-TextValue "D28" "29.47"
$ws.Range("E28").Value = "  -5.01%  "
Set-TextValue "D29" "7.67"
$ws.Range("E29").Value = "  -3.32%  "
Set-TextValue "D30" "0.118"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("E33").Value = "  -0.09%  "
Set-TextValue "D34" "41.12"
$ws.Range("E34").Value = "  -2.20%  "
Set-TextValue "D35" "2.49"
$ws.Range("E35").Value = "  -2.63%  "
$ws.Range("E36").Value = "  -0.02%  "
Set-TextValue "D37" "52.13"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("E38").Value = "  +0.01%  "
Set-TextValue "D39" "3.43"
$ws.Range("E39").Value = "  -2.22%  "
Set-TextValue "D40" "2.93"
$ws.Range("E40").Value = "  -2.76%  "
Set-TextValue "D41" "137.89"
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  +5.72%  "
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("E45").Value = "  +5.64%  "
Set-TextValue "D46" "16.69"
$ws.Range("E46").Value = "  -2.48%  "
$ws.Range("E47").Value = "  +1.35%  "
Set-TextValue "D48" "21.35"
$ws.Range("E48").Value = "  -2.52%  "
Set-TextValue "D49" "2.130.81"
$ws.Range("E49").Value = "  -3.42%  "
$ws.Range("E50").Value = "  -4.82%  "
$ws.Range("E51").Value = "  +0.38%  "
